$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33,8).Value = 642.13043
$ws.Cells.Item(33,9).Value = 496.64285
$ws.Cells.Item(33,10).Value = 868.44446
$ws.Cells.Item(33,11).Value = 496.64285
$ws.Cells.Item(33,12).Value = 868.44446
$ws.Cells.Item(33,13).Value = -267.64285
$ws.Cells.Item(33,14).Value = -1326.44446
$ws.Cells.Item(112,8).Value = 7813488
$ws.Cells.Item(112,9).Value = 400
$ws.Cells.Item(112,10).Value = 9616508
$ws.Cells.Item(112,11).Value = 1200
$ws.Cells.Item(112,12).Value = 28849524
$ws.Cells.Item(112,13).Value = -92
$ws.Cells.Item(112,14).Value = -28851740
$ws.Cells.Item(137,8).Value = 5268137
$ws.Cells.Item(137,9).Value = 8339475
$ws.Cells.Item(137,10).Value = 2985.7144
$ws.Cells.Item(137,11).Value = 25018425
$ws.Cells.Item(137,12).Value = 8957.143199999999
$ws.Cells.Item(137,13).Value = -25015875
$ws.Cells.Item(137,14).Value = -14057.1432
$ws.Cells.Item(141,8).Value = 829007.2
$ws.Cells.Item(141,9).Value = 1691.76
$ws.Cells.Item(141,10).Value = 4276155
$ws.Cells.Item(141,11).Value = 5075.28
$ws.Cells.Item(141,12).Value = 12828465
$ws.Cells.Item(141,13).Value = 104.7200000000003
$ws.Cells.Item(141,14).Value = -12838825
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4,8).Value = 149.66667
$ws.Cells.Item(4,9).Value = 99.333336
$ws.Cells.Item(4,10).Value = 200
$ws.Cells.Item(4,11).Value = 99.333336
$ws.Cells.Item(4,12).Value = 200
$ws.Cells.Item(4,13).Value = 16.666664
$ws.Cells.Item(4,14).Value = -432
$ws.Cells.Item(132,8).Value = 29415902
$ws.Cells.Item(132,9).Value = 41670796
$ws.Cells.Item(132,10).Value = 4152.6
$ws.Cells.Item(132,11).Value = 125012388
$ws.Cells.Item(132,12).Value = 12457.8
$ws.Cells.Item(132,13).Value = -125009858
$ws.Cells.Item(132,14).Value = -17517.8
$ws.Cells.Item(133,8).Value = 29992
$ws.Cells.Item(133,10).Value = 29992
$ws.Cells.Item(133,12).Value = 29992
$ws.Cells.Item(133,14).Value = -35052
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64,8).Value = 452.63635
$ws.Cells.Item(64,9).Value = 355.14285
$ws.Cells.Item(64,10).Value = 498.13333
$ws.Cells.Item(64,11).Value = 355.14285
$ws.Cells.Item(64,12).Value = 498.13333
$ws.Cells.Item(64,13).Value = -130.14285
$ws.Cells.Item(64,14).Value = -948.13333
$ws.Cells.Item(67,8).Value = 452.63635
$ws.Cells.Item(67,9).Value = 355.14285
$ws.Cells.Item(67,10).Value = 498.13333
$ws.Cells.Item(67,11).Value = 355.14285
$ws.Cells.Item(67,12).Value = 498.13333
$ws.Cells.Item(67,13).Value = 424.85715
$ws.Cells.Item(67,14).Value = -2058.13333
$ws.Cells.Item(86,8).Value = 1744.4762
$ws.Cells.Item(86,9).Value = 1196.25
$ws.Cells.Item(86,11).Value = 1196.25
$ws.Cells.Item(86,13).Value = -73.25
$ws.Cells.Item(89,8).Value = 1744.4762
$ws.Cells.Item(89,9).Value = 1196.25
$ws.Cells.Item(89,11).Value = 5981.25
$ws.Cells.Item(89,13).Value = -365.25
$ws.Cells.Item(105,8).Value = 1432.2258
$ws.Cells.Item(105,9).Value = 1220.6522
$ws.Cells.Item(105,10).Value = 2040.5
$ws.Cells.Item(105,11).Value = 1220.6522
$ws.Cells.Item(105,12).Value = 2040.5
$ws.Cells.Item(105,13).Value = 526.3478
$ws.Cells.Item(105,14).Value = -5534.5
$ws.Cells.Item(134,8).Value = 2298.5173
$ws.Cells.Item(134,9).Value = 1415
$ws.Cells.Item(134,10).Value = 5075.2856
$ws.Cells.Item(134,11).Value = 4245
$ws.Cells.Item(134,12).Value = 15225.8568
$ws.Cells.Item(134,13).Value = -1710
$ws.Cells.Item(134,14).Value = -20295.8568
$ws.Cells.Item(135,8).Value = 31183.334
$ws.Cells.Item(135,10).Value = 31183.334
$ws.Cells.Item(135,12).Value = 31183.334
$ws.Cells.Item(135,14).Value = -41323.334
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16,8).Value = 5199.8
$ws.Cells.Item(16,9).Value = 8999
$ws.Cells.Item(16,10).Value = 4250
$ws.Cells.Item(16,11).Value = 8999
$ws.Cells.Item(16,12).Value = 4250
$ws.Cells.Item(16,13).Value = -8712
$ws.Cells.Item(16,14).Value = -4824
$ws.Cells.Item(25,8).Value = 30829.75
$ws.Cells.Item(25,10).Value = 34376.855
$ws.Cells.Item(25,12).Value = 34376.855
$ws.Cells.Item(25,14).Value = -34724.855
$ws.Cells.Item(31,8).Value = 2086492
$ws.Cells.Item(31,9).Value = 2779873
$ws.Cells.Item(31,10).Value = 6349.1665
$ws.Cells.Item(31,11).Value = 2779873
$ws.Cells.Item(31,12).Value = 6349.1665
$ws.Cells.Item(31,13).Value = -2779578
$ws.Cells.Item(31,14).Value = -6939.1665
$ws.Cells.Item(34,8).Value = 2086492
$ws.Cells.Item(34,9).Value = 2779873
$ws.Cells.Item(34,10).Value = 6349.1665
$ws.Cells.Item(34,11).Value = 2779873
$ws.Cells.Item(34,12).Value = 6349.1665
$ws.Cells.Item(34,13).Value = -2779671
$ws.Cells.Item(34,14).Value = -6753.1665
$ws.Cells.Item(58,8).Value = 20002938
$ws.Cells.Item(58,9).Value = 1523.8182
$ws.Cells.Item(58,10).Value = 35718336
$ws.Cells.Item(58,11).Value = 1523.8182
$ws.Cells.Item(58,12).Value = 35718336
$ws.Cells.Item(58,13).Value = -1320.8182
$ws.Cells.Item(58,14).Value = -35718742
$ws.Cells.Item(113,8).Value = 5199.8
$ws.Cells.Item(113,9).Value = 8999
$ws.Cells.Item(113,10).Value = 4250
$ws.Cells.Item(113,11).Value = 8999
$ws.Cells.Item(113,12).Value = 4250
$ws.Cells.Item(113,13).Value = -6829
$ws.Cells.Item(113,14).Value = -8590
$ws.Cells.Item(134,8).Value = 3763.3
$ws.Cells.Item(134,9).Value = 1740
$ws.Cells.Item(134,10).Value = 5786.6
$ws.Cells.Item(134,11).Value = 5220
$ws.Cells.Item(134,12).Value = 17359.8
$ws.Cells.Item(134,13).Value = -2685
$ws.Cells.Item(134,14).Value = -22429.8
$ws.Cells.Item(136,8).Value = 20002938
$ws.Cells.Item(136,9).Value = 1523.8182
$ws.Cells.Item(136,10).Value = 35718336
$ws.Cells.Item(136,11).Value = 4571.4546
$ws.Cells.Item(136,12).Value = 107155008
$ws.Cells.Item(136,13).Value = -2021.4546
$ws.Cells.Item(136,14).Value = -107160108
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98,8).Value = 198.5
$ws.Cells.Item(98,10).Value = 175
$ws.Cells.Item(98,12).Value = 525
$ws.Cells.Item(98,14).Value = -3521
$ws.Cells.Item(107,8).Value = 901.92
$ws.Cells.Item(107,10).Value = 833.41174
$ws.Cells.Item(107,12).Value = 2500.23522
$ws.Cells.Item(107,14).Value = -6340.23522
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21,8).Value = 53338
$ws.Cells.Item(21,9).Value = 20000
$ws.Cells.Item(21,10).Value = 70007
$ws.Cells.Item(21,11).Value = 20000
$ws.Cells.Item(21,12).Value = 70007
$ws.Cells.Item(21,14).Value = -70353
$ws.Cells.Item(21,13).Value = -19827
$ws.Cells.Item(24,8).Value = 627164
$ws.Cells.Item(24,9).Value = 1001935.3
$ws.Cells.Item(24,10).Value = 65007
$ws.Cells.Item(24,11).Value = 1001935.3
$ws.Cells.Item(24,12).Value = 65007
$ws.Cells.Item(24,14).Value = -65353
$ws.Cells.Item(24,13).Value = -1001762.3
$ws.Cells.Item(30,8).Value = 53338
$ws.Cells.Item(30,9).Value = 20000
$ws.Cells.Item(30,10).Value = 70007
$ws.Cells.Item(30,11).Value = 20000
$ws.Cells.Item(30,12).Value = 70007
$ws.Cells.Item(30,14).Value = -70217
$ws.Cells.Item(30,13).Value = -19895
$ws.Cells.Item(132,8).Value = 2968.0322
$ws.Cells.Item(132,9).Value = 2600.647
$ws.Cells.Item(132,10).Value = 3414.1428
$ws.Cells.Item(132,11).Value = 7801.941
$ws.Cells.Item(132,12).Value = 10242.4284
$ws.Cells.Item(132,13).Value = -5271.941
$ws.Cells.Item(132,14).Value = -15302.4284
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2,8).Value = 287714.56
$ws.Cells.Item(2,9).Value = 324412.06
$ws.Cells.Item(2,10).Value = 131750.25
$ws.Cells.Item(2,11).Value = 324412.06
$ws.Cells.Item(2,12).Value = 131750.25
$ws.Cells.Item(2,13).Value = -324300.06
$ws.Cells.Item(2,14).Value = -131974.25
$ws.Cells.Item(7,8).Value = 3200.6667
$ws.Cells.Item(7,9).Value = 1602
$ws.Cells.Item(7,10).Value = 4000
$ws.Cells.Item(7,11).Value = 1602
$ws.Cells.Item(7,12).Value = 4000
$ws.Cells.Item(7,13).Value = -1490
$ws.Cells.Item(7,14).Value = -4224
$ws.Cells.Item(14,8).Value = 645253.25
$ws.Cells.Item(14,9).Value = 837002.7
$ws.Cells.Item(14,10).Value = 70005
$ws.Cells.Item(14,11).Value = 837002.7
$ws.Cells.Item(14,12).Value = 70005
$ws.Cells.Item(14,13).Value = -836830.7
$ws.Cells.Item(14,14).Value = -70349
$ws.Cells.Item(16,8).Value = 732.3333
$ws.Cells.Item(16,9).Value = 849.75
$ws.Cells.Item(16,10).Value = 497.5
$ws.Cells.Item(16,11).Value = 849.75
$ws.Cells.Item(16,12).Value = 497.5
$ws.Cells.Item(16,13).Value = -679.75
$ws.Cells.Item(16,14).Value = -837.5
$ws.Cells.Item(24,8).Value = 37503.5
$ws.Cells.Item(24,9).Value = 0
$ws.Cells.Item(24,10).Value = 37503.5
$ws.Cells.Item(24,11).Value = 0
$ws.Cells.Item(24,12).Value = 37503.5
$ws.Cells.Item(24,14).Value = -38189.5
$ws.Cells.Item(24,13).ClearContents()
$ws.Cells.Item(68,8).Value = 1486.9546
$ws.Cells.Item(68,10).Value = 3588.25
$ws.Cells.Item(68,12).Value = 3588.25
$ws.Cells.Item(68,14).Value = -5086.25
$ws.Cells.Item(71,8).Value = 1486.9546
$ws.Cells.Item(71,10).Value = 3588.25
$ws.Cells.Item(71,12).Value = 17941.25
$ws.Cells.Item(71,14).Value = -25429.25
$ws.Cells.Item(100,8).Value = 2850
$ws.Cells.Item(100,9).Value = 1450
$ws.Cells.Item(100,10).Value = 3550
$ws.Cells.Item(100,11).Value = 1450
$ws.Cells.Item(100,12).Value = 3550
$ws.Cells.Item(100,13).Value = -909
$ws.Cells.Item(100,14).Value = -4632
$ws.Cells.Item(126,8).Value = 3200.6667
$ws.Cells.Item(126,9).Value = 1602
$ws.Cells.Item(126,10).Value = 4000
$ws.Cells.Item(126,11).Value = 4806
$ws.Cells.Item(126,12).Value = 12000
$ws.Cells.Item(126,13).Value = -2336
$ws.Cells.Item(126,14).Value = -16940
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2,8).Value = 3520001
$ws.Cells.Item(2,9).Value = 1006666.3
$ws.Cells.Item(2,10).Value = 5028001.5
$ws.Cells.Item(2,11).Value = 1006666.3
$ws.Cells.Item(2,12).Value = 5028001.5
$ws.Cells.Item(2,13).Value = -1006554.3
$ws.Cells.Item(2,14).Value = -5028225.5
$ws.Cells.Item(107,8).Value = 2455.5557
$ws.Cells.Item(107,9).Value = 1657.1428
$ws.Cells.Item(107,10).Value = 5250
$ws.Cells.Item(107,11).Value = 4971.428400000001
$ws.Cells.Item(107,12).Value = 15750
$ws.Cells.Item(107,13).Value = -3051.428400000001
$ws.Cells.Item(107,14).Value = -19590
$ws.Cells.Item(135,8).Value = 68803.75
$ws.Cells.Item(135,10).Value = 68803.75
$ws.Cells.Item(135,12).Value = 68803.75
$ws.Cells.Item(135,14).Value = -78943.75
$ws.Cells.Item(136,8).Value = 2350.8
$ws.Cells.Item(136,9).Value = 1536.2
$ws.Cells.Item(136,10).Value = 3980
$ws.Cells.Item(136,11).Value = 4608.6
$ws.Cells.Item(136,12).Value = 11940
$ws.Cells.Item(136,13).Value = -2058.6
$ws.Cells.Item(136,14).Value = -17040
